# Add the missed name "Dane" to the attendance list, and move the
# "_GoBack" bookmark (Word auto-tracks the last edit location) from the
# end of the document to right after the newly inserted text.

$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark wherever Word currently has it
# (Word recreates this automatically at the last edit point, so we
# delete the stale one before inserting new text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert ", Dane" right after "Christian Hidalgo" in the attendance line.
$found = $d.Content.Find.Execute("Christian Hidalgo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Christian Hidalgo, Dane", 2)

# Re-add the _GoBack bookmark at the point right after the inserted text,
# matching Word's behaviour of marking the last edited location.
$again = $d.Content.Find.Execute("Christian Hidalgo, Dane", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$r = $d.Content.Duplicate
$r.Start = $d.Content.End
$r.End = $d.Content.End

$found2 = $d.Content.Find.Execute("Christian Hidalgo, Dane", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
if ($found2) {
    $endRange = $d.Content.Find.Parent
}

$d.Bookmarks.Add("_GoBack", $word.Selection.Range) | Out-Null
